# Block Company Client,Brand,Product Changes
#
# Reorders the "Block*" test-case rows (2-9) on the "IND_Critical Regression"
# sheet, flips most of their "Execute" flags (and a few further down the
# sheet) to "No", and moves the active tab / selection from
# "CHN_SIT_FullCycle" to "IND_Critical Regression" (cell D7:D8 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IND_Critical Regression")

# --- Reorder / rewrite the "Block*" rows (UnitName / TestCases / Description / Execute) ---
$blockRows = @(
    @(2, "BlockCompanyBrand",  "Blockcompanybrand",  "Blockcompanybrand",            "Yes"),
    @(3, "BlockCompanyProduct","Blockcompanyproduct","Blockcompanyproduct",          "No"),
    @(4, "BlockCompanyClient", "BlockCompanyclient", "BlockCompanyclient",           "No"),
    @(5, "BlockGlobalBrand",   "Blockglobalbrand",   "Block created global brand",   "No"),
    @(6, "BlockGlobalProduct", "Blockglobalproduct", "Blockglobalproduct",           "No"),
    @(7, "BlockGlobalClient",  "Blockglobalclient",  "Block created global client",  "No"),
    @(8, "BlockCompanyVendor", "BlockCompanyVendor", "BlockCompanyVendor",           "No"),
    @(9, "BlockGlobalVendor",  "BlockGlobalVendor",  "Block Created Global Vendor",  "No")
)

foreach ($row in $blockRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# --- Flip the remaining "Execute" flags further down the sheet to "No" ---
$ws.Range("D30").Value = "No"
$ws.Range("D35").Value = "No"
$ws.Range("D36").Value = "No"
$ws.Range("D37").Value = "No"
$ws.Range("D38").Value = "No"

# --- Move the active tab to "IND_Critical Regression" with D7:D8 selected ---
$ws.Activate()
$ws.Range("D7:D8").Select()
